# Update Name of Algo
# Applies updated KNN imputation result values to specific cells in Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "D7"   = -7.453
    "A8"   = -22.188
    "A10"  = -21.846
    "A12"  = -21.589
    "D15"  = -8.147
    "A18"  = -22.269
    "D18"  = -8.366
    "D20"  = -7.562
    "D29"  = -7.292
    "D30"  = -7.179
    "D31"  = -8.023
    "A37"  = -20.21700000000001
    "D40"  = -7.662999999999999
    "D50"  = -8.104999999999999
    "A55"  = -22.283
    "A68"  = -21.534
    "D68"  = -6.778
    "D76"  = -7.753
    "A77"  = -20.934
    "A78"  = -20.281
    "A81"  = -21.811
    "A82"  = -22.261
    "D87"  = -8.261999999999999
    "D88"  = -8.293000000000001
    "D96"  = -7.267
    "D98"  = -8.242000000000001
    "D101" = -7.986999999999999
    "D102" = -8.142999999999999
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
